$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- AA2: drop the cell border, keep the quote-prefix formatting ---------
$ws.Range("AA2").Borders.LineStyle = -4142   # xlLineStyleNone

# --- AA3: give it the same look as AC2/AC3 (quote-prefixed, Arial 10, #444444, bordered)
# then fill in the new value "EMEAAD\pvergez" (added as a new shared string) ---
$ws.Range("AC2").Copy()
$ws.Range("AA3").PasteSpecial(-4122)          # xlPasteFormats
$ws.Range("AA3").Formula = "'EMEAAD\pvergez"

# --- update the view: move the selection to AB3 --------------------------
$null = $ws.Range("AB3").Select()
